$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for 2021-03-15 ("Ají", variety
# "Inferno", "Primera" quality). It belongs right above the existing
# row 36 (chronologically it sorts there), so insert a fresh row at
# position 36 and shift every following record down by one — this is
# exactly what the XML diff shows (old row N -> new row N+1, for every
# N from 36 to 91, plus one brand-new row 36 and the dimension growing
# from R91 to R92).
$ws.Rows.Item(36).Insert()

$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 44771
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112021
$ws.Range("G36").Value = "Ají"
$ws.Range("H36").Value = "Inferno"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 120
$ws.Range("K36").Value = 11000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 11500
$ws.Range("N36").Value = "$/caja 15 kilos"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 767
$ws.Range("Q36").Value = 15
$ws.Range("R36").Value = "Hortaliza"
